# PlayerPerformance_3924.xlsx - add "ODI Bowling Extra" scraping sheet,
# and tidy up placeholder blank cells left over on "ODI Batting Extra".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Trim the stray empty placeholder cells on "ODI Batting Extra" (sheet 4).
#    Rows 2-7, 11, 12 keep column A (MATCH_CODE) and F (MAN_OF_MATCH) but
#    lose the empty B/C/D/E placeholders. Row 8 only had a blank C/D/E (its
#    B8 already holds a real MAIDEN value). Rows 9-10 already hold real data
#    in every column and are left untouched. Rows 13-21 only ever had empty
#    placeholders beyond column A, so every one of B:F is cleared there.
# ---------------------------------------------------------------------------
$battingExtra = $wb.Worksheets.Item("ODI Batting Extra")
$battingExtra.Range("B2:E7").ClearContents()
$battingExtra.Range("C8:E8").ClearContents()
$battingExtra.Range("B11:E12").ClearContents()
$battingExtra.Range("B13:F21").ClearContents()

# ---------------------------------------------------------------------------
# 2. Add the new "ODI Bowling Extra" sheet after "ODI Batting Extra".
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$bowlingExtra = $wb.Worksheets.Add($null, $lastSheet)
$bowlingExtra.Name = "ODI Bowling Extra"

# Header row
$bowlingExtra.Cells.Item(1, 1).Value = "MATCH_CODE"
$bowlingExtra.Cells.Item(1, 2).Value = "MAIDEN_OVERS"
$bowlingExtra.Cells.Item(1, 3).Value = "PERCENT_WICKETS_OF_ALL"
$bowlingExtra.Range("A1:C1").Font.Bold = $true

# Every value on this sheet is stored as text in the source export, so force
# the whole used range to Text format before writing any values - this stops
# the COM layer from auto-coercing "0"/"10.00%" into numbers/percentages.
$dataRange = $bowlingExtra.Range("A2:C21")
$dataRange.NumberFormat = "@"

$rows = @(
    @("4305", "0", "10.00%"),
    @("4311", "", ""),
    @("4315", "0", ""),
    @("4328", "", ""),
    @("4333", "", ""),
    @("4337", "", ""),
    @("4341", "1", "40.00%"),
    @("4346", "", ""),
    @("4353", "", ""),
    @("4355", "", ""),
    @("4423", "0", ""),
    @("4452", "0", "40.00%"),
    @("4453", "0", "10.00%"),
    @("4455", "1", ""),
    @("4636", "0", "20.00%"),
    @("4639", "1", "30.00%"),
    @("4642", "", ""),
    @("4647", "", ""),
    @("4648", "2", "40.00%"),
    @("4649", "4", "20.00%")
)

$r = 2
foreach ($row in $rows) {
    $bowlingExtra.Cells.Item($r, 1).Value = $row[0]
    $bowlingExtra.Cells.Item($r, 2).Value = $row[1]
    $bowlingExtra.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
